$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)
$target = $para1.Characters(27, 11)
$target.Text = "first few times "
